$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$shape = $s.Shapes.Item(3)
$table = $shape.Table
$table.ApplyStyle("{AEF2AB11-503C-4B30-86FA-78EF3D482B21}")
